$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Name / account holder info
$ws.Range("C2").Value = "Hartmut"

# Card number must remain text (16-digit string). Force a text number
# format so COM doesn't coerce the digit string to a numeric value, write
# it, then restore the original (General) cell style by pasting the
# formatting from a neighbouring cell that still carries that style - this
# avoids leaving behind a stray "text" number-format style on the cell.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("D3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Value = "Mohaupt"

# Opening balance line
$ws.Range("D5").Value = "KONTOSTAND AM 29.03.2024"

# Transaction row 6
$ws.Range("B6").Value = "31.03."
$ws.Range("C6").Value = "01.04."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-3630567"
$ws.Range("E6").Value = "53,45-"

# Transaction row 7
$ws.Range("B7").Value = "01.04."
$ws.Range("C7").Value = "02.04."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 1656036"
$ws.Range("E7").Value = "41,77-"

# Transaction row 8
$ws.Range("B8").Value = "05.04."
$ws.Range("C8").Value = "06.04."
$ws.Range("D8").Value = "MCDONALDS Cottbus"
$ws.Range("E8").Value = "14,05-"

# Transaction row 9
$ws.Range("B9").Value = "06.04."
$ws.Range("C9").Value = "07.04."
$ws.Range("D9").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E9").Value = "54,20-"

# Closing balance line
$ws.Range("D12").Value = "KONTOSTAND AM 09.04.2024"
$ws.Range("E12").Value = "163,47-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 15.04.2024"
